$d = $word.ActiveDocument

# Locate the paragraph that contains "Current Recommendation" so we can
# rebuild it with the new run structure described by the diff.
$rng = $d.Content.Duplicate
$found = $rng.Find.Execute("Current Recommendation")
if (-not $found) {
    throw "Could not find the 'Current Recommendation' paragraph"
}

$p = $rng.Paragraphs(1)
$pRange = $p.Range

# Rebuild the paragraph's OOXML: the bold "Current Recommendation" run is
# kept as-is, and the long run that used to hold the whole explanation is
# split into the nine runs described by the diff.
$newParaXml = '<w:p w14:paraId="1E0B6036" w14:textId="77777777" w:rsidR="00A43E6E" w:rsidRDefault="00000000">' +
  '<w:r><w:rPr><w:b/></w:rPr><w:t>Current Recommendation</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">: As these are 3 of the 8 most important features according to our model training, we have reason for Big Mountain Resort to increase their price. </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">We recommend Big Mountain Resort to up their price by $13. </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">With an average of 350,000 visitors purchasing an average of 5 tickets each season, </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">this increase in ticket </w:t></w:r>' +
  '<w:r><w:t>price would increase revenue by $</w:t></w:r>' +
  '<w:r><w:t>22,</w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">750,000 </w:t></w:r>' +
  '<w:r><w:t xml:space="preserve">well above </w:t></w:r>' +
  '<w:r><w:t>the 1,540,000 dollars needed for the new chair lift.</w:t></w:r>' +
  '</w:p>'

$payload = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
  '<pkg:xmlData>' +
  '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:body>' + $newParaXml + '</w:body>' +
  '</w:document>' +
  '</pkg:xmlData>' +
  '</pkg:part>' +
  '</pkg:package>'

$null = $pRange.InsertXML($payload)
